# Auto-generated edit script: updates numeric ticket-count (F) and price (G) cells
# across all 4 worksheets to match the target workbook snapshot.
$wb = $excel.ActiveWorkbook

# --- 展览 (Worksheets.Item(1)) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1305
$ws.Range("F4").Value = 171
$ws.Range("F5").Value = 962
$ws.Range("G5").Value = 70
$ws.Range("F6").Value = 29
$ws.Range("F7").Value = 1048
$ws.Range("F8").Value = 164
$ws.Range("G8").Value = 78
$ws.Range("F9").Value = 847
$ws.Range("G9").Value = 68.8
$ws.Range("F10").Value = 51
$ws.Range("F11").Value = 728
$ws.Range("G11").Value = 59
$ws.Range("F12").Value = 1417
$ws.Range("F13").Value = 1038
$ws.Range("F14").Value = 759
$ws.Range("F15").Value = 770
$ws.Range("F16").Value = 89
$ws.Range("F17").Value = 592
$ws.Range("F18").Value = 105
$ws.Range("F19").Value = 664
$ws.Range("F20").Value = 1274
$ws.Range("F21").Value = 185
$ws.Range("F22").Value = 180
$ws.Range("F23").Value = 5282
$ws.Range("F24").Value = 276
$ws.Range("F26").Value = 2449
$ws.Range("F27").Value = 5889
$ws.Range("F29").Value = 1005
$ws.Range("F30").Value = 597
$ws.Range("F31").Value = 66
$ws.Range("F33").Value = 1051
$ws.Range("F35").Value = 43
$ws.Range("F37").Value = 696
$ws.Range("F39").Value = 50
$ws.Range("F41").Value = 1081
$ws.Range("F44").Value = 9
$ws.Range("F45").Value = 26
$ws.Range("F46").Value = 94
$ws.Range("F47").Value = 494

# --- 演出 (Worksheets.Item(2)) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 10
$ws.Range("F5").Value = 2088
$ws.Range("F6").Value = 74
$ws.Range("F8").Value = 122
$ws.Range("G8").Value = 108
$ws.Range("F9").Value = 483
$ws.Range("F11").Value = 102
$ws.Range("F12").Value = 96
$ws.Range("F13").Value = 134
$ws.Range("F15").Value = 673
$ws.Range("F17").Value = 786
$ws.Range("F18").Value = 17
$ws.Range("F28").Value = 14
$ws.Range("F31").Value = 152
$ws.Range("F35").Value = 69
$ws.Range("F37").Value = 15
$ws.Range("F38").Value = 97
$ws.Range("F39").Value = 902
$ws.Range("F40").Value = 494
$ws.Range("F42").Value = 29
$ws.Range("F43").Value = 2
$ws.Range("F45").Value = 79
$ws.Range("F46").Value = 99
$ws.Range("F48").Value = 14

# --- 本地生活 (Worksheets.Item(3)) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 668
$ws.Range("F5").Value = 764
$ws.Range("F6").Value = 415
$ws.Range("F7").Value = 238

# --- 全部类型 (Worksheets.Item(4)) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 668
$ws.Range("F4").Value = 1305
$ws.Range("F5").Value = 764
$ws.Range("F7").Value = 415
$ws.Range("F8").Value = 238
$ws.Range("F9").Value = 238
$ws.Range("F10").Value = 74
$ws.Range("F11").Value = 171
$ws.Range("F12").Value = 962
$ws.Range("G12").Value = 70
$ws.Range("F14").Value = 1048
$ws.Range("F15").Value = 164
$ws.Range("G15").Value = 78
$ws.Range("F16").Value = 847
$ws.Range("G16").Value = 68.8
$ws.Range("F17").Value = 51
$ws.Range("F18").Value = 727
$ws.Range("G18").Value = 59
$ws.Range("F19").Value = 1417
$ws.Range("F20").Value = 102
$ws.Range("F21").Value = 96
$ws.Range("F22").Value = 1038
$ws.Range("F23").Value = 759
$ws.Range("F24").Value = 134
$ws.Range("F25").Value = 770
$ws.Range("F26").Value = 1274
$ws.Range("F27").Value = 185
$ws.Range("F28").Value = 17
$ws.Range("F30").Value = 180
$ws.Range("F32").Value = 276
$ws.Range("F34").Value = 2449
$ws.Range("F35").Value = 5889
$ws.Range("F37").Value = 66
$ws.Range("F39").Value = 43
$ws.Range("F40").Value = 696
$ws.Range("F41").Value = 50
$ws.Range("F42").Value = 1081
$ws.Range("F44").Value = 902
$ws.Range("F45").Value = 494
$ws.Range("F46").Value = 29
$ws.Range("F47").Value = 94
$ws.Range("F49").Value = 79
$ws.Range("F50").Value = 14
